# "unit info in progress" -- add the Solapur generating-unit rows, a new
# eminating_lines column, and switch trial_run_date to a real date format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Gadarwara) values: ramp 800 -> 75, and the
# trial_run_date becomes an actual date (2019-06-03) instead of "1600".
$ws.Range("C2").Value = 75
$ws.Range("D2").Value = 43619

# Add new rows 3 and 4 for the Solapur units (1 and 2). "Solapur" is written
# before "eminating_lines" below so it lands earlier in the shared-strings
# table, matching the saved workbook.
$ws.Range("A3").Value = "Solapur"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = 43618

$ws.Range("A4").Value = "Solapur"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 75
$ws.Range("D4").Value = 43618

# Now add the new "eminating_lines" column with its data
$ws.Range("E1").Value = "eminating_lines"
$ws.Range("E2").Value = 8
$ws.Range("E3").Value = 6
$ws.Range("E4").Value = 5

# Apply date number format to column D (entire column, so the column-level
# style is recorded just like column C's wrap-text style)
$ws.Columns("D").NumberFormat = "[$-14009]yyyy/mm/dd;@"

# Set the selection to mirror the saved view state
$ws.Range("E5").Select()

# Configure the page setup (paper size / orientation) for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
